$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates from the Oct 16 2023 GitHub Actions crypto-price refresh.
# NumberFormat "@" + ClearFormats() keeps numeric-looking strings (e.g. "213.70",
# "0.0600") stored as text instead of being coerced to doubles, while leaving the
# cell style index untouched (matches the unstyled cells in the source workbook).
$updates = @{
    "D2" = "28.413.87"
    "E2" = "  +4.16%  "
    "D3" = "1.582.86"
    "E3" = "  +0.34%  "
    "D4" = "0.998"
    "E4" = "  -1.16%  "
    "D5" = "213.70"
    "E5" = "  +0.87%  "
    "D6" = "0.497"
    "E6" = "  +0.51%  "
    "D7" = "0.998"
    "E7" = "  -1.04%  "
    "D8" = "23.80"
    "E8" = "  +7.74%  "
    "E9" = "  +0.82%  "
    "D10" = "0.0600"
    "D11" = "0.0886"
    "E11" = "  +1.84%  "
    "D12" = "1.808.26"
    "E12" = "  +0.42%  "
    "D13" = "1.589.14"
    "E13" = "  +1.04%  "
    "D14" = "3.77"
    "E14" = "  -0.74%  "
    "E15" = "  +1.38%  "
    "D16" = "28.348.71"
    "E16" = "  +4.13%  "
    "D17" = "63.93"
    "E17" = "  +2.42%  "
    "D18" = "233.07"
    "E18" = "  +7.45%  "
    "E19" = "  +0.63%  "
    "D20" = "7.47"
    "E20" = "  -0.02%  "
    "E21" = "  -1.03%  "
    "D22" = "4.13"
    "E22" = "  -0.73%  "
    "D23" = "9.37"
    "E23" = "  +1.07%  "
    "E24" = "  -0.46%  "
    "D25" = "151.48"
    "E25" = "  -1.77%  "
    "D26" = "15.30"
    "E26" = "  +1.11%  "
    "E27" = "  -1.22%  "
    "E28" = "  -0.16%  "
    "E29" = "  -0.82%  "
    "E30" = "  -0.09%  "
    "E31" = "  -0.14%  "
    "E32" = "  -0.56%  "
    "E33" = "  -1.31%  "
    "D34" = "1.416.20"
    "E34" = "  -2.81%  "
    "E35" = "  -1.44%  "
    "E36" = "  -5.79%  "
    "E37" = "  -1.63%  "
    "E38" = "  -0.35%  "
    "E39" = "  +7.45%  "
    "D40" = "0.544"
    "E40" = "  +1.24%  "
    "D41" = "0.812"
    "E41" = "  -0.22%  "
    "E42" = "  -1.10%  "
    "E44" = "  +5.01%  "
    "D45" = "0.972"
    "E45" = "  -3.47%  "
    "D46" = "64.32"
    "E46" = "  -0.73%  "
    "D47" = "1.718.42"
    "E47" = "  +0.59%  "
    "D48" = "87.27"
    "E48" = "  +1.49%  "
    "B49" = "BabyDogeCoin"
    "C49" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D49" = "0.0₆0103"
    "E49" = "  +0.97%  "
    "B50" = "Cronos"
    "C50" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D50" = "0.0525"
    "E50" = "  +0.04%  "
    "D51" = "39.46"
    "E51" = "  +16.18%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}

Write-Host "Applied $($updates.Count) cell updates"
